$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 531 (old "44407" week),
# shifting all subsequent rows down by 2. This grows the used range from
# A1:R550 to A1:R552, matching the new weekly data being added.
$ws.Rows.Item(531).Insert()
$ws.Rows.Item(531).Insert()

# Populate the two newly inserted rows (531 = "Primera", 532 = "Segunda")
# with the new week's values. The non-changing descriptive columns mirror
# every other row in this block (same market / region / product / unit).

# Row 531 - Calidad "Primera"
$ws.Range("A531").Value = 8
$ws.Range("B531").Value = "Terminal La Palmera de La Serena"
$ws.Range("C531").Value = "Coquimbo"
$ws.Range("D531").Value = 44939
$ws.Range("E531").Value = 4
$ws.Range("F531").Value = 100112017
$ws.Range("G531").Value = "Apio"
$ws.Range("H531").Value = "Americana (o)"
$ws.Range("I531").Value = "Primera"
$ws.Range("J531").Value = 1700
$ws.Range("K531").Value = 9000
$ws.Range("L531").Value = 10000
$ws.Range("M531").Value = 9500
$ws.Range("N531").Value = "`$/docena de matas"
$ws.Range("O531").Value = "Provincia del Elquí"
$ws.Range("P531").Value = 1583
$ws.Range("Q531").Value = 6
$ws.Range("R531").Value = "Hortaliza"

# Row 532 - Calidad "Segunda"
$ws.Range("A532").Value = 8
$ws.Range("B532").Value = "Terminal La Palmera de La Serena"
$ws.Range("C532").Value = "Coquimbo"
$ws.Range("D532").Value = 44939
$ws.Range("E532").Value = 4
$ws.Range("F532").Value = 100112017
$ws.Range("G532").Value = "Apio"
$ws.Range("H532").Value = "Americana (o)"
$ws.Range("I532").Value = "Segunda"
$ws.Range("J532").Value = 900
$ws.Range("K532").Value = 7000
$ws.Range("L532").Value = 8000
$ws.Range("M532").Value = 7500
$ws.Range("N532").Value = "`$/docena de matas"
$ws.Range("O532").Value = "Provincia del Elquí"
$ws.Range("P532").Value = 1250
$ws.Range("Q532").Value = 6
$ws.Range("R532").Value = "Hortaliza"

# Match the date-format style already used by the other "Fecha" cells.
$ws.Range("D531").NumberFormat = $ws.Range("D533").NumberFormat
$ws.Range("D532").NumberFormat = $ws.Range("D533").NumberFormat
